# Updated symbol list (coin rankings) per commit "Updated symbol list on
# Thu Dec 15 18:57:03 UTC 2022 with GitHub Actions" -- refreshes Price
# column values and re-ranks several rows whose coins changed position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column (D) values are stored as text in this workbook (e.g. "260.70"),
# so force a Text number format before assigning, then restore the default
# 'Normal' style so no stray formatting is left behind.
function Set-TextCell($sheet, $ref, $value) {
    $cell = $sheet.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" '260.09'
Set-TextCell $ws "D3" '22.72'
Set-TextCell $ws "D4" '6.215'
Set-TextCell $ws "D5" '0.06076'
Set-TextCell $ws "D6" '3.516'
Set-TextCell $ws "D7" '6.716'
Set-TextCell $ws "D9" '0.8000'
Set-TextCell $ws "D10" '0.1577'
Set-TextCell $ws "D11" '0.08091'
Set-TextCell $ws "D12" '0.03323'
Set-TextCell $ws "D13" '0.03124'
Set-TextCell $ws "D14" '0.09270'
Set-TextCell $ws "D15" '3.911'
Set-TextCell $ws "D16" '0.001708'
Set-TextCell $ws "D17" '0.04827'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextCell $ws "D18" '0.0006172'
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("B19").Value = 'TigerCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell $ws "D19" '0.006206'
$ws.Range("E19").Value = '18TigerCashTCH'
$ws.Range("B20").Value = 'BitKan'
$ws.Range("C20").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextCell $ws "D20" '0.001103'
$ws.Range("E20").Value = '19BitKanKAN'
$ws.Range("B21").Value = 'HotbitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextCell $ws "D21" '0.003376'
$ws.Range("E21").Value = '20HotbitTokenHTB'
$ws.Range("B22").Value = 'NitroEx'
$ws.Range("C22").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextCell $ws "D22" '0.0001502'
$ws.Range("E22").Value = '21NitroExNTX'
$ws.Range("B23").Value = 'LEO'
$ws.Range("C23").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell $ws "D23" '3.691'
$ws.Range("E23").Value = '22LEOLEO'
$ws.Range("B24").Value = 'BTSEToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextCell $ws "D24" '2.294'
$ws.Range("E24").Value = '23BTSETokenBTSE'
Set-TextCell $ws "D25" '0.3373'
Set-TextCell $ws "D27" '0.0006179'
Set-TextCell $ws "D40" '0.04602'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextCell $ws "D41" '0.1119'
$ws.Range("E41").Value = '40BKEXTokenBKK'
Set-TextCell $ws "D42" '0.003134'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextCell $ws "D43" '0.003398'
$ws.Range("E43").Value = '42KickTokenKICK'
Set-TextCell $ws "D45" '0.002976'
Set-TextCell $ws "D46" '0.00006022'
Set-TextCell $ws "D47" '0.00000000752'
Set-TextCell $ws "D48" '0.7516'
Set-TextCell $ws "D49" '0.1319'
Set-TextCell $ws "D50" '0.00001503'
Set-TextCell $ws "D51" '0.01012'
